$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$addr, [string]$val) {
    # Force text interpretation (avoids Excel auto-converting
    # numeric-looking strings like "244.31" into numbers) without
    # touching NumberFormat/style: build the literal via a TEXT
    # formula, then paste-special VALUES ONLY into the target cell.
    $ws.Range("ZZ1").Formula = '="' + $val + '"'
    $ws.Range("ZZ1").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

$ws.Range("D2").Value = "36.715.06"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "2.060.81"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue "D5" "244.31"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue "D8" "55.45"
$ws.Range("E8").Value = "  -5.69%  "
Set-TextValue "D9" "59.76"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("E10").Value = "  -3.46%  "
Set-TextValue "D11" "0.0754"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("E13").Value = "  +5.53%  "
Set-TextValue "D14" "14.87"
$ws.Range("E14").Value = "  -3.55%  "
$ws.Range("D15").Value = "2.360.23"
$ws.Range("E15").Value = "  +0.35%  "
Set-TextValue "D16" "5.48"
$ws.Range("E16").Value = "  -4.76%  "
$ws.Range("D17").Value = "2.062.03"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "36.628.58"
$ws.Range("E18").Value = "  -1.07%  "
Set-TextValue "D19" "17.20"
$ws.Range("E19").Value = "  -5.49%  "
Set-TextValue "D20" "72.24"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").Value = "0.0₃0868"
$ws.Range("E21").Value = "  -2.15%  "
Set-TextValue "D22" "238.81"
$ws.Range("E22").Value = "  +0.24%  "
Set-TextValue "D23" "5.29"
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -2.41%  "
Set-TextValue "D26" "2.18"
$ws.Range("E26").Value = "  +1.51%  "
Set-TextValue "D27" "9.36"
$ws.Range("E27").Value = "  -4.23%  "
Set-TextValue "D28" "165.39"
$ws.Range("E28").Value = "  -2.08%  "
Set-TextValue "D29" "20.24"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("E30").Value = "  -1.39%  "
Set-TextValue "D31" "5.14"
$ws.Range("E31").Value = "  -7.92%  "
$ws.Range("E32").Value = "  +6.37%  "
$ws.Range("E33").Value = "  -3.96%  "
Set-TextValue "D34" "0.0601"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("E35").Value = "  +0.11%  "
Set-TextValue "D36" "1.84"
$ws.Range("E36").Value = "  -0.35%  "
Set-TextValue "D37" "0.0849"
$ws.Range("E37").Value = "  -0.15%  "
Set-TextValue "D38" "2.23"
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("E40").Value = "  -4.48%  "
$ws.Range("E41").Value = "  -4.89%  "
$ws.Range("E42").Value = "  -2.85%  "
$ws.Range("E43").Value = "  -3.54%  "
Set-TextValue "D44" "95.09"
$ws.Range("E44").Value = "  -2.89%  "
Set-TextValue "D45" "7.75"
$ws.Range("E45").Value = "  +14.74%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D46" "0.0906"
$ws.Range("E46").Value = "  -6.94%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.411.23"
$ws.Range("E47").Value = "  +8.36%  "
Set-TextValue "D48" "16.15"
$ws.Range("E48").Value = "  -5.03%  "
$ws.Range("E49").Value = "  +1.43%  "
Set-TextValue "D50" "2.29"
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("D51").Value = "2.249.56"
$ws.Range("E51").Value = "  +0.68%  "

$ws.Range("ZZ1").Clear()
